# Rename the existing sheet to "TestCases" and add a new "TestData" sheet after it
# (copy, rather than Add, so the new sheet inherits the same x14ac namespace setup).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TestCases"

$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TestData"

# ---- TestCases sheet content ----
$ws1.Range("A1").Value = "TestCases"
$ws1.Range("B1").Value = "RunMode"
$ws1.Range("A2").Value = "AddCustomerTest"
$ws1.Range("A3").Value = "OpenAccountTest"
$ws1.Range("B2").Value = "y"
$ws1.Range("B3").Value = "n"

$ws1.Columns.Item(1).ColumnWidth = 18.6
$ws1.Columns.Item(2).ColumnWidth = 14.6

# ---- TestData sheet content ----
# AddCustomerTest block
$ws2.Range("A1").Value = "AddCustomerTest"
$ws2.Range("A2").Value = "runMode"
$ws2.Range("B2").Value = "firstName"
$ws2.Range("C2").Value = "lastName"
$ws2.Range("D2").Value = "postalCode"
$ws2.Range("A3").Value = "y"
$ws2.Range("B3").Value = "Serge"
$ws2.Range("C3").Value = "Koko"
$ws2.Range("D3").Value = 12345
$ws2.Range("A4").Value = "n"
$ws2.Range("B4").Value = "John"
$ws2.Range("C4").Value = "Doe"
$ws2.Range("D4").Value = 98125

# OpenAccountTest block
$ws2.Range("A6").Value = "OpenAccountTest"
$ws2.Range("A7").Value = "runMode"
$ws2.Range("B7").Value = "customer"
$ws2.Range("C7").Value = "currency"
$ws2.Range("A8").Value = "y"
$ws2.Range("C8").Value = "Dollar"
$ws2.Range("B8").Value = "Serge Koko"
$ws2.Range("A9").Value = "y"
$ws2.Range("B9").Value = "Harry Potter"
$ws2.Range("C9").Value = "Dollar"
$ws2.Range("A10").Value = "n"
$ws2.Range("B10").Value = "John Doe"
$ws2.Range("C10").Value = "EURO"

$ws2.Columns.Item(1).ColumnWidth = 13.25
$ws2.Columns.Item(2).ColumnWidth = 13.25
$ws2.Columns.Item(3).ColumnWidth = 13.25
$ws2.Columns.Item(4).ColumnWidth = 13.25

# ---- Selections ----
$ws2.Range("C10").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B10").Select() | Out-Null
